# Daily TGP (terminal gate pricing) update
# Rolls forward effective dates and updates price columns (Diesel/ULP/PULP/e10)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value2 = 45959
$ws.Range("D8").Value2 = 164.26
$ws.Range("E8").Value2 = 159.51
$ws.Range("F8").Value2 = 169.51
$ws.Range("G8").Value2 = 159.67

$ws.Range("A9").Value2 = 45959
$ws.Range("D9").Value2 = 164.26
$ws.Range("E9").Value2 = 159.51
$ws.Range("F9").Value2 = 169.51
$ws.Range("G9").Value2 = 159.67

$ws.Range("A10").Value2 = 45959
$ws.Range("D10").Value2 = 167.11
$ws.Range("E10").Value2 = 162.22
$ws.Range("F10").Value2 = 172.22
$ws.Range("G10").Value2 = 162.69

$ws.Range("A11").Value2 = 45958
$ws.Range("D11").Value2 = 163.75
$ws.Range("E11").Value2 = 159.53
$ws.Range("F11").Value2 = 169.53
$ws.Range("G11").Value2 = 159.69

$ws.Range("A12").Value2 = 45958
$ws.Range("D12").Value2 = 163.75
$ws.Range("E12").Value2 = 159.53
$ws.Range("F12").Value2 = 169.53
$ws.Range("G12").Value2 = 159.69

$ws.Range("A13").Value2 = 45958
$ws.Range("D13").Value2 = 166.18
$ws.Range("E13").Value2 = 161.81
$ws.Range("F13").Value2 = 171.81
$ws.Range("G13").Value2 = 162.28

$ws.Range("A17").Value2 = 45959
$ws.Range("D17").Value2 = 170.56
$ws.Range("E17").Value2 = 165.08
$ws.Range("F17").Value2 = 175.08

$ws.Range("A18").Value2 = 45958
$ws.Range("D18").Value2 = 169.63
$ws.Range("E18").Value2 = 164.59
$ws.Range("F18").Value2 = 174.59

$ws.Range("A22").Value2 = 45959
$ws.Range("D22").Value2 = 166.05
$ws.Range("E22").Value2 = 161.23
$ws.Range("F22").Value2 = 170.83
$ws.Range("G22").Value2 = 162.4

$ws.Range("A23").Value2 = 45959
$ws.Range("D23").Value2 = 171.88
$ws.Range("E23").Value2 = 165.94
$ws.Range("F23").Value2 = 175.94

$ws.Range("A24").Value2 = 45959
$ws.Range("D24").Value2 = 171.69
$ws.Range("E24").Value2 = 166.14
$ws.Range("F24").Value2 = 176.14

$ws.Range("A25").Value2 = 45959
$ws.Range("D25").Value2 = 172.51
$ws.Range("E25").Value2 = 165.53
$ws.Range("F25").Value2 = 175.53
$ws.Range("G25").Value2 = 165.36

$ws.Range("A26").Value2 = 45959
$ws.Range("D26").Value2 = 171.24
$ws.Range("E26").Value2 = 167.08
$ws.Range("F26").Value2 = 177.08

$ws.Range("A27").Value2 = 45958
$ws.Range("D27").Value2 = 165.11
$ws.Range("E27").Value2 = 160.81
$ws.Range("F27").Value2 = 170.41
$ws.Range("G27").Value2 = 161.99

$ws.Range("A28").Value2 = 45958
$ws.Range("D28").Value2 = 170.95
$ws.Range("E28").Value2 = 165.53
$ws.Range("F28").Value2 = 175.53

$ws.Range("A29").Value2 = 45958
$ws.Range("D29").Value2 = 170.76
$ws.Range("E29").Value2 = 165.73
$ws.Range("F29").Value2 = 175.73

$ws.Range("A30").Value2 = 45958
$ws.Range("D30").Value2 = 171.59
$ws.Range("E30").Value2 = 165.12
$ws.Range("F30").Value2 = 175.12
$ws.Range("G30").Value2 = 164.95

$ws.Range("A31").Value2 = 45958
$ws.Range("D31").Value2 = 170.31
$ws.Range("E31").Value2 = 166.68
$ws.Range("F31").Value2 = 176.68

$ws.Range("A35").Value2 = 45959
$ws.Range("D35").Value2 = 165.35
$ws.Range("E35").Value2 = 159.43
$ws.Range("F35").Value2 = 168.43

$ws.Range("A36").Value2 = 45958
$ws.Range("D36").Value2 = 164.43
$ws.Range("E36").Value2 = 159.02
$ws.Range("F36").Value2 = 168.02

$ws.Range("A40").Value2 = 45959
$ws.Range("D40").Value2 = 170.99
$ws.Range("E40").Value2 = 164.77
$ws.Range("F40").Value2 = 174.77

$ws.Range("A41").Value2 = 45959
$ws.Range("D41").Value2 = 170.7
$ws.Range("E41").Value2 = 165.19
$ws.Range("F41").Value2 = 175.19

$ws.Range("A42").Value2 = 45958
$ws.Range("D42").Value2 = 170.08
$ws.Range("E42").Value2 = 164.3
$ws.Range("F42").Value2 = 174.3

$ws.Range("A43").Value2 = 45958
$ws.Range("D43").Value2 = 169.79
$ws.Range("E43").Value2 = 164.72
$ws.Range("F43").Value2 = 174.72

$ws.Range("A47").Value2 = 45959
$ws.Range("D47").Value2 = 163.64
$ws.Range("E47").Value2 = 160.38
$ws.Range("F47").Value2 = 170.38

$ws.Range("A48").Value2 = 45959
$ws.Range("D48").Value2 = 163.63
$ws.Range("E48").Value2 = 160.55
$ws.Range("F48").Value2 = 170.55

$ws.Range("A49").Value2 = 45958
$ws.Range("D49").Value2 = 161.88
$ws.Range("E49").Value2 = 159.86
$ws.Range("F49").Value2 = 169.86

$ws.Range("A50").Value2 = 45958
$ws.Range("D50").Value2 = 161.87
$ws.Range("E50").Value2 = 160.04
$ws.Range("F50").Value2 = 170.04

$ws.Range("A54").Value2 = 45959
$ws.Range("D54").Value2 = 181.2
$ws.Range("E54").Value2 = 175.11
$ws.Range("F54").Value2 = 185.11

$ws.Range("A55").Value2 = 45959
$ws.Range("D55").Value2 = 168.85
$ws.Range("E55").Value2 = 172.62
$ws.Range("F55").Value2 = 182.62

$ws.Range("A56").Value2 = 45959
$ws.Range("D56").Value2 = 171.12

$ws.Range("A57").Value2 = 45959
$ws.Range("D57").Value2 = 170.82
$ws.Range("E57").Value2 = 166.89

$ws.Range("A58").Value2 = 45959
$ws.Range("D58").Value2 = 166.72
$ws.Range("E58").Value2 = 162.94
$ws.Range("F58").Value2 = 172.94

$ws.Range("A59").Value2 = 45959
$ws.Range("D59").Value2 = 173.55
$ws.Range("E59").Value2 = 173.34

$ws.Range("A60").Value2 = 45958
$ws.Range("D60").Value2 = 180.28
$ws.Range("E60").Value2 = 174.75
$ws.Range("F60").Value2 = 184.75

$ws.Range("A61").Value2 = 45958
$ws.Range("D61").Value2 = 167.93
$ws.Range("E61").Value2 = 171.98
$ws.Range("F61").Value2 = 181.98

$ws.Range("A62").Value2 = 45958
$ws.Range("D62").Value2 = 170.2

$ws.Range("A63").Value2 = 45958
$ws.Range("D63").Value2 = 169.88
$ws.Range("E63").Value2 = 166.25

$ws.Range("A64").Value2 = 45958
$ws.Range("D64").Value2 = 165.78
$ws.Range("E64").Value2 = 162.3
$ws.Range("F64").Value2 = 172.3

$ws.Range("A65").Value2 = 45958
$ws.Range("D65").Value2 = 172.61
$ws.Range("E65").Value2 = 172.96
